$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 1759410.1
$ws.Range("I92").Value = 3078452.8
$ws.Range("K92").Value = 3078452.8
$ws.Range("M92").Value = -3077204.8
$ws.Range("H103").Value = 887.46875
$ws.Range("J103").Value = 800.0769
$ws.Range("L103").Value = 2400.2307
$ws.Range("N103").Value = -3572.2307
$ws.Range("H137").Value = 2030.2307
$ws.Range("I137").Value = 1223.3334
$ws.Range("K137").Value = 3670.0002
$ws.Range("M137").Value = -1120.0002

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").Value = ""
$ws.Range("H47").Value = 63041
$ws.Range("J47").Value = 63041
$ws.Range("L47").Value = 63041
$ws.Range("N47").Value = -64491
$ws.Range("H61").Value = 5029.1
$ws.Range("I61").Value = 2706.5715
$ws.Range("J61").Value = 10448.333
$ws.Range("K61").Value = 2706.5715
$ws.Range("L61").Value = 10448.333
$ws.Range("M61").Value = -2494.5715
$ws.Range("N61").Value = -10872.333
$ws.Range("H74").Value = 1077.4375
$ws.Range("I74").Value = 1002.7143
$ws.Range("J74").Value = 1600.5
$ws.Range("K74").Value = 1002.7143
$ws.Range("L74").Value = 1600.5
$ws.Range("M74").Value = -128.7143
$ws.Range("N74").Value = -3348.5
$ws.Range("H77").Value = 1077.4375
$ws.Range("I77").Value = 1002.7143
$ws.Range("J77").Value = 1600.5
$ws.Range("K77").Value = 5013.5715
$ws.Range("L77").Value = 8002.5
$ws.Range("M77").Value = -645.5715
$ws.Range("N77").Value = -16738.5
$ws.Range("H132").Value = 1515.6666
$ws.Range("I132").Value = 1212.7778
$ws.Range("K132").Value = 3638.3334
$ws.Range("M132").Value = -1108.3334
$ws.Range("H136").Value = 5029.1
$ws.Range("I136").Value = 2706.5715
$ws.Range("J136").Value = 10448.333
$ws.Range("K136").Value = 8119.7145
$ws.Range("L136").Value = 31344.999
$ws.Range("M136").Value = -5569.7145
$ws.Range("N136").Value = -36444.999

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1516.5
$ws.Range("I99").Value = 1100
$ws.Range("K99").Value = 1100
$ws.Range("M99").Value = 398
$ws.Range("H105").Value = 1709.5186
$ws.Range("I105").Value = 1736.8462
$ws.Range("K105").Value = 1736.8462
$ws.Range("M105").Value = 10.15380000000005
$ws.Range("H134").Value = 17425.533
$ws.Range("I134").Value = 20274.908
$ws.Range("J134").Value = 9589.75
$ws.Range("K134").Value = 60824.724
$ws.Range("L134").Value = 28769.25
$ws.Range("M134").Value = -58289.724
$ws.Range("N134").Value = -33839.25

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2717.3333
$ws.Range("I31").Value = 1843.7391
$ws.Range("J31").Value = 5587.7144
$ws.Range("K31").Value = 1843.7391
$ws.Range("L31").Value = 5587.7144
$ws.Range("M31").Value = -1548.7391
$ws.Range("N31").Value = -6177.7144
$ws.Range("H34").Value = 2717.3333
$ws.Range("I34").Value = 1843.7391
$ws.Range("J34").Value = 5587.7144
$ws.Range("K34").Value = 1843.7391
$ws.Range("L34").Value = 5587.7144
$ws.Range("M34").Value = -1641.7391
$ws.Range("N34").Value = -5991.7144
$ws.Range("H58").Value = 1403838.5
$ws.Range("I58").Value = 2718469.8
$ws.Range("K58").Value = 2718469.8
$ws.Range("M58").Value = -2718266.8
$ws.Range("H132").Value = 2085.087
$ws.Range("I132").Value = 1063.6111
$ws.Range("K132").Value = 3190.8333
$ws.Range("M132").Value = -660.8333000000002
$ws.Range("H134").Value = 1848.5555
$ws.Range("I134").Value = 1369.1177
$ws.Range("J134").Value = 9999
$ws.Range("K134").Value = 4107.3531
$ws.Range("L134").Value = 29997
$ws.Range("M134").Value = -1572.3531
$ws.Range("N134").Value = -35067
$ws.Range("H136").Value = 1403838.5
$ws.Range("I136").Value = 2718469.8
$ws.Range("K136").Value = 8155409.399999999
$ws.Range("M136").Value = -8152859.399999999

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 2999.5
$ws.Range("I22").Value = 3000
$ws.Range("J22").Value = 2999
$ws.Range("K22").Value = 9000
$ws.Range("L22").Value = 8997
$ws.Range("M22").Value = -8831
$ws.Range("N22").Value = -9335
$ws.Range("H27").Value = 2999.5
$ws.Range("I27").Value = 3000
$ws.Range("J27").Value = 2999
$ws.Range("K27").Value = 9000
$ws.Range("L27").Value = 8997
$ws.Range("M27").Value = -8898
$ws.Range("N27").Value = -9201
$ws.Range("H39").Value = 1366.3334
$ws.Range("J39").Value = 1366.3334
$ws.Range("L39").Value = 4099.0002
$ws.Range("N39").Value = -4687.0002
$ws.Range("H43").Value = 4000
$ws.Range("J43").Value = 4000
$ws.Range("L43").Value = 12000
$ws.Range("N43").Value = -12228
$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").Value = ""
$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").Value = ""
$ws.Range("H131").Value = 779.1429000000001
$ws.Range("J131").Value = 798.9674
$ws.Range("L131").Value = 2396.9022
$ws.Range("N131").Value = -12476.9022
$ws.Range("H132").Value = 1710.5
$ws.Range("I132").Value = 1650
$ws.Range("J132").Value = 1730.6666
$ws.Range("K132").Value = 14850
$ws.Range("L132").Value = 15575.9994
$ws.Range("M132").Value = -12320
$ws.Range("N132").Value = -20635.9994

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 15000
$ws.Range("J43").Value = 15000
$ws.Range("L43").Value = 15000
$ws.Range("N43").Value = -15302
$ws.Range("H97").Value = 1270.625
$ws.Range("I97").Value = 1357.2222
$ws.Range("J97").Value = 1159.2858
$ws.Range("K97").Value = 1357.2222
$ws.Range("L97").Value = 1159.2858
$ws.Range("M97").Value = -861.2221999999999
$ws.Range("N97").Value = -2151.2858
$ws.Range("H126").Value = 2461344.8
$ws.Range("I126").Value = 3706133
$ws.Range("K126").Value = 11118399
$ws.Range("M126").Value = -11115929
$ws.Range("H132").Value = 1927650.9
$ws.Range("I132").Value = 3499425.2
$ws.Range("K132").Value = 10498275.6
$ws.Range("M132").Value = -10495745.6

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7170
$ws.Range("I7").Value = 3525.75
$ws.Range("K7").Value = 3525.75
$ws.Range("M7").Value = -3413.75
$ws.Range("H126").Value = 7170
$ws.Range("I126").Value = 3525.75
$ws.Range("K126").Value = 10577.25
$ws.Range("M126").Value = -8107.25
$ws.Range("H132").Value = 2311.7144
$ws.Range("I132").Value = 2199
$ws.Range("J132").Value = 2356.8
$ws.Range("K132").Value = 6597
$ws.Range("L132").Value = 7070.400000000001
$ws.Range("M132").Value = -4067
$ws.Range("N132").Value = -12130.4
$ws.Range("H134").Value = 48507.332
$ws.Range("J134").Value = 48507.332
$ws.Range("L134").Value = 48507.332
$ws.Range("N134").Value = -58647.332
$ws.Range("H136").Value = 4144.154
$ws.Range("I136").Value = 1863.1666
$ws.Range("J136").Value = 6099.2856
$ws.Range("K136").Value = 5589.4998
$ws.Range("L136").Value = 18297.8568
$ws.Range("M136").Value = -3039.4998
$ws.Range("N136").Value = -23397.8568

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 7532.933
$ws.Range("I132").Value = 2499.25
$ws.Range("J132").Value = 9363.362999999999
$ws.Range("K132").Value = 7497.75
$ws.Range("L132").Value = 28090.089
$ws.Range("M132").Value = -4967.75
$ws.Range("N132").Value = -33150.089
$ws.Range("H136").Value = 39684972
$ws.Range("I136").Value = 92594200
$ws.Range("K136").Value = 277782600
$ws.Range("M136").Value = -277780050
